$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 84
$ws.Range("I2").Value = 204
$ws.Range("J2").Value = 789
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 218
$ws.Range("M2").Value = 11
$ws.Range("N2").Value = 145
$ws.Range("Q2").Value = 1
$ws.Range("R2").Value = 13
$ws.Range("S2").Value = 86
$ws.Range("T2").Value = 147
$ws.Range("V2").Value = 1149
$ws.Range("W2").Value = 0
$ws.Range("X2").Value = 1140
$ws.Range("Y2").Value = 4
$ws.Range("Z2").Value = 18
$ws.Range("AA2").Value = 5
